# Update "想去人数" (want-to-go count) values in column F
# for the "展览" sheet and the "全部类型" sheet.
# Each entry maps a row number to its new F value.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 58
    3  = 118
    4  = 1669
    5  = 297
    7  = 2039
    8  = 10299
    11 = 264
    14 = 7161
    15 = 1106
    17 = 102
    19 = 251
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 58
    3  = 118
    4  = 1669
    5  = 297
    8  = 2039
    11 = 10299
    14 = 264
    17 = 7161
    18 = 1106
    20 = 102
    22 = 251
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
